$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A9").Value = 9
$ws.Range("B9").Value = 'url="www.myapplication.com"'

$ws.Range("G10").Select()
